# Refresh the cryptos price table with the latest scraped values.
# (GitHub Actions run: "Updated cryptos list" data refresh.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.256.91"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "2.950.09"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'568.11"
$ws.Range("E5").Value = "  -2.94%  "
$ws.Range("D6").Value = "'159.31"
$ws.Range("E6").Value = "  +3.13%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.518"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "2.946.96"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").Value = "'6.65"
$ws.Range("E10").Value = "  -4.98%  "
$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("E12").Value = "  +2.05%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "'34.27"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "65.368.56"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "3.440.80"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").Value = "'6.97"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "2.959.53"
$ws.Range("E19").Value = "  -1.37%  "
$ws.Range("D20").Value = "'14.91"
$ws.Range("E20").Value = "  +8.36%  "
$ws.Range("D21").Value = "'445.91"
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("D22").Value = "'0.687"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "'7.24"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("D24").Value = "'82.20"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("E25").Value = "  -1.94%  "
$ws.Range("D26").Value = "'12.08"
$ws.Range("E26").Value = "  -3.90%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "'10.00"
$ws.Range("E28").Value = "  -6.97%  "
$ws.Range("E29").Value = "  +2.01%  "
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").Value = "'2.57"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("D33").Value = "'27.17"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").Value = "'0.978"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("D37").Value = "'5.75"
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").Value = "'48.99"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("B39").Value = "Arweave"
$ws.Range("C39").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D39").Value = "'44.25"
$ws.Range("E39").Value = "  -3.64%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.98"
$ws.Range("E40").Value = "  -8.37%  "
$ws.Range("E41").Value = "  -1.59%  "
$ws.Range("D42").Value = "'2.83"
$ws.Range("E42").Value = "  -3.78%  "
$ws.Range("D43").Value = "'0.297"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").Value = "'385.30"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "2.709.99"
$ws.Range("E47").Value = "  -2.17%  "
$ws.Range("D48").Value = "'132.97"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D50").Value = "'2.19"
$ws.Range("E50").Value = "  +5.20%  "
$ws.Range("E51").Value = "  +0.19%  "
